$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.850.30"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.082.88"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.15"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.37"
$ws.Range("E7").Value = "  +3.67%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.75"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.21"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.077.54"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.756.01"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.17"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.39"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.138"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.99"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0217"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.73"
$ws.Range("E43").Value = "  +7.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.445.18"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.39"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.274.40"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.90"
$ws.Range("E51").Value = "  +1.07%  "
